# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45178 (2023-09-09) to 45179 (2023-09-10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45178) {
        $cell.Value2 = 45179
    }
}
